$d = $word.ActiveDocument

# Locate the paragraph that currently ends with the unfinished sentence
# "... fait combiner les deux faces car ils enlevent les aretes entre les deux face"
# and the (bookmark-only) paragraph right after it.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "fait combiner les deux faces") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

$firstPara = $d.Paragraphs.Item($targetIndex)
$secondPara = $d.Paragraphs.Item($targetIndex + 1)

# Range spanning the whole first paragraph through the end (incl. the
# paragraph mark) of the paragraph right after it - both get rewritten as a
# single merged paragraph followed by a brand new, fully empty paragraph.
$r = $d.Range($firstPara.Range.Start, $secondPara.Range.End)

$xml = '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Malheureusement, ça ne marche pas souvent comme prévu, et il faut à tâtons, à l’aide des autres icônes de cet onglet, essayer de réparer ces saletés de face</w:t></w:r><w:r><w:t>s, en testant les outils possibles de réparations</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:br/><w:t>Souvent, l’outil « </w:t></w:r><w:r><w:t>enlever les arêtes surnuméraires</w:t></w:r><w:r><w:t> »</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>fait combiner les deux faces car ils enlèvent les arêtes entre les deux face</w:t></w:r><w:r><w:t xml:space="preserve">s. Néanmoins, quand on rouvre le fichier STEP, les arêtes surnuméraires ne sont pas </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>enlevés</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">. Mais chaque </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>modif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> des</w:t></w:r><w:r><w:t xml:space="preserve"> autres</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> icones de la partie Fix/Réparer et </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Adjust</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/Ajuster sont conservés, donc il faut absolument les faire.</w:t></w:r></w:p><w:p/>'

[void]$r.InsertXML($xml)
